{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Two changes (per the commit \"Modify outro; remove messages\"):\n// 1. Remove the whole \"## Adding missing grouping variables: `sample_id`\"\n//    source-code paragraph (a leftover console/message block).\n// 2. Shorten/rewrite the closing sentences of the final \"Outro\" paragraph.\n\nconst body = context.document.body;\n\n// --- 1. Remove the stray \"Adding missing grouping variables\" paragraph ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst stray = paragraphs.items.find((p) =>\n  p.text.indexOf(\"## Adding missing grouping variables: `sample_id`\") !== -1\n);\nif (stray) {\n  stray.delete();\n}\nawait context.sync();\n\n// --- 2. Rewrite the outro text ---\nconst oldOutro =\n  \"uses some of its functionality. It would be great to further develop \" +\n  \"streamlining packages which are good at doing single things. It would \" +\n  \"also be fantastic to co-develop a new set of programs that automatically \" +\n  \"tune spectral machine learning pipelines. Complex problems require \" +\n  \"targeted learning. For example, one could create a custom graph learner \" +\n  \"using mlr3 and a preprocessing wrapper targeted to spectral analysis, in \" +\n  \"connection with a proper database system. If you have ideas, just send \" +\n  \"me an email or interact via github.\";\n\nconst newOutro =\n  \"uses some of its functionality. Complex problems and professional \" +\n  \"spectroscopy applications require transfer learning and spectral \" +\n  \"feature engineering pipelines that tune automatically. If you have \" +\n  \"ideas to collaborate and develop new frameworks, just send me an email \" +\n  \"or interact via github.\";\n\nconst hits = body.search(oldOutro, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(newOutro, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Two changes (per the commit \"Modify outro; remove messages\"):\n# 1. Remove the whole \"## Adding missing grouping variables: `sample_id`\"\n#    source-code paragraph (a leftover console/message block).\n# 2. Shorten/rewrite the closing sentences of the final \"Outro\" paragraph.\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the stray \"Adding missing grouping variables\" paragraph ---\n$strayMarker = \"## Adding missing grouping variables: ``sample_id``\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($strayMarker)) {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n\n# --- 2. Rewrite the outro text ---\n$oldOutro = \"uses some of its functionality. It would be great to further develop streamlining packages which are good at doing single things. It would also be fantastic to co-develop a new set of programs that automatically tune spectral machine learning pipelines. Complex problems require targeted learning. For example, one could create a custom graph learner using mlr3 and a preprocessing wrapper targeted to spectral analysis, in connection with a proper database system. If you have ideas, just send me an email or interact via github.\"\n$newOutro = \"uses some of its functionality. Complex problems and professional spectroscopy applications require transfer learning and spectral feature engineering pipelines that tune automatically. If you have ideas to collaborate and develop new frameworks, just send me an email or interact via github.\"\n\n$rng = $d.Content\n$rng.Find.Text = $oldOutro\n$found = $rng.Find.Execute()\nif ($found) {\n    $rng.Text = $newOutro\n}\n"}
